# Refresh the "cryptos" price/volume table (Price = column D, Volume(1h) = column E).
# Note: some Price values (e.g. '213.28') look like plain numbers to Excel and
# would otherwise be auto-converted from text to a numeric cell; a leading
# apostrophe forces them to stay text, matching the source data, while Excel
# strips the apostrophe itself from the stored value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.158.74'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '1.612.36'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = '''213.28'
$ws.Range('E5').Value = '  +1.99%  '
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('D7').Value = '''0.483'
$ws.Range('E7').Value = '  +1.12%  '
$ws.Range('E8').Value = '  +2.24%  '
$ws.Range('D9').Value = '''0.0620'
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').Value = '''18.48'
$ws.Range('E10').Value = '  +3.81%  '
$ws.Range('D11').Value = '''0.0796'
$ws.Range('E11').Value = '  +1.27%  '
$ws.Range('D12').Value = '1.831.75'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').Value = '1.612.01'
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('D14').Value = '''4.07'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').Value = '''0.513'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').Value = '26.140.97'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '''60.93'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  +3.70%  '
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').Value = '''199.71'
$ws.Range('E20').Value = '  +5.58%  '
$ws.Range('D21').Value = '''4.28'
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').Value = '''9.50'
$ws.Range('E22').Value = '  +2.42%  '
$ws.Range('D23').Value = '''6.03'
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('E24').Value = '  +3.49%  '
$ws.Range('D25').Value = '''142.93'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').Value = '''1.73'
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').Value = '''15.27'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').Value = '''6.53'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('E30').Value = '  -2.10%  '
$ws.Range('E31').Value = '  +3.06%  '
$ws.Range('E32').Value = '  +3.12%  '
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('D36').Value = '1.106.72'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.0153'
$ws.Range('E37').Value = '  +1.77%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''0.509'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''2.34'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  -0.65%  '
$ws.Range('D41').Value = '''0.790'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').Value = '''0.801'
$ws.Range('E42').Value = '  +8.38%  '
$ws.Range('D43').Value = '''5.18'
$ws.Range('E43').Value = '  +2.37%  '
$ws.Range('D44').Value = '1.743.70'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').Value = '''93.32'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '''1.54'
$ws.Range('E47').Value = '  +6.36%  '
$ws.Range('D48').Value = '''54.08'
$ws.Range('E48').Value = '  +2.14%  '
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').Value = '''1.00'
$ws.Range('E51').Value = '  -0.55%  '
